# DPLKAKT067-001 - Setup Periode Bulanan - "Update Regresi Tanggal 31/03/2023"
# Roll the test-data period forward one year: 2023 -> 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TGL_AWAL (O2): 15/04/2023 -> 15/04/2024 (real date cell)
$ws.Range("O2").Value = [DateTime]"2024-04-15"

# TGL_AKHIR (P2): 15/05/2023 -> 15/05/2024 (stored as text, matches source layout)
$ws.Range("P2").Value = "15/05/2024"

# PERIODE_BULANAN (Q2): 202305 -> 202405
$ws.Range("Q2").Value = "202405"

# Move the viewport / selection the way the author left it: scrolled right to
# column N and the active cell on W2.
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1

$ws.Range("W2").Select()
